# Fix bug in fedrollover: extend the date series in column A / B down to
# cover the rest of August 2020 (rows 373-379), matching the existing
# pattern of incrementing the "day" portion of the pseudo-date code by
# 100 each row, with a value of 0 in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 373
$startValue = 20082500

for ($i = 0; $i -lt 7; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startValue + ($i * 100)
    $ws.Cells.Item($row, 2).Value = 0
}
